# Delete rows 19 and 20 (the "125v1genav37-black" product image entries),
# shifting all subsequent rows up by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A19:F20").EntireRow.Delete() | Out-Null
